$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# Sheet1 (Schedule) changes: 22
$ws1.Range("B2").Value = 46057.25
$ws1.Range("C2").Value = 6
$ws1.Range("D2").Value = 22.68
$ws1.Range("E2").Value = 1075.50601275
$ws1.Range("F2").Value = 47.42090003306879
$ws1.Range("A3").Value = 46057.27083333334
$ws1.Range("B3").Value = 46057.66666666666
$ws1.Range("C3").Value = 9.5
$ws1.Range("D3").Value = 35.91
$ws1.Range("E3").Value = 370.5886860000001
$ws1.Range("F3").Value = 10.31992999164578
$ws1.Range("A4").Value = 46057.97916666666
$ws1.Range("B4").Value = 46058.14583333334
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 15.12
$ws1.Range("E4").Value = 1102.29294825
$ws1.Range("F4").Value = 72.90297276785715
$ws1.Range("A5").Value = 46058.29166666666
$ws1.Range("B5").Value = 46058.64583333334
$ws1.Range("C5").Value = 8.5
$ws1.Range("D5").Value = 32.13
$ws1.Range("E5").Value = 657.68489475
$ws1.Range("F5").Value = 20.46949563492064
# Row 5 is brand new; match date/time display format used by rows 2-4 in columns A:B
$ws1.Range("A5:B5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Sheet2 (Detailed) changes: 68
$ws2.Range("E14").Value = "OFF"
$ws2.Range("B38").Value = 66.67267
$ws2.Range("B39").Value = 101.75293
$ws2.Range("B40").Value = 154.65275
$ws2.Range("C40").Value = "historical"
$ws2.Range("B41").Value = 142.86493
$ws2.Range("C41").Value = "historical"
$ws2.Range("B42").Value = 169.25335
$ws2.Range("C42").Value = "historical"
$ws2.Range("B43").Value = 158.53104
$ws2.Range("C43").Value = "historical"
$ws2.Range("B44").Value = 156.16131
$ws2.Range("C44").Value = "historical"
$ws2.Range("B45").Value = 133.73657
$ws2.Range("C45").Value = "historical"
$ws2.Range("B46").Value = 115.03999
$ws2.Range("C46").Value = "historical"
$ws2.Range("B47").Value = 151.9658
$ws2.Range("C47").Value = "historical"
$ws2.Range("B48").Value = 138.42
$ws2.Range("C48").Value = "historical"
$ws2.Range("B49").Value = 161.47181
$ws2.Range("E49").Value = "ON"
$ws2.Range("B50").Value = 142.5809
$ws2.Range("E50").Value = "ON"
$ws2.Range("B51").Value = 138.42
$ws2.Range("E51").Value = "ON"
$ws2.Range("B52").Value = 138.42
$ws2.Range("B54").Value = 136.85635
$ws2.Range("B55").Value = 136.27451
$ws2.Range("B56").Value = 138.1133
$ws2.Range("B57").Value = 138.22945
$ws2.Range("E57").Value = "OFF"
$ws2.Range("B58").Value = 137.32372
$ws2.Range("E58").Value = "OFF"
$ws2.Range("E59").Value = "OFF"
$ws2.Range("B60").Value = 146.32223
$ws2.Range("B61").Value = 161.73439
$ws2.Range("B62").Value = 161.43949
$ws2.Range("B63").Value = 150.83269
$ws2.Range("B64").Value = 96.13648999999999
$ws2.Range("E64").Value = "ON"
$ws2.Range("B65").Value = 57.08
$ws2.Range("B66").Value = 53.94454
$ws2.Range("B68").Value = 24.73456
$ws2.Range("B69").Value = 24.7421
$ws2.Range("B70").Value = 25.30266
$ws2.Range("B73").Value = 9.221959999999999
$ws2.Range("B74").Value = 20.57901
$ws2.Range("B77").Value = 36.0601
$ws2.Range("B78").Value = 36.0601
$ws2.Range("B79").Value = 36.0601
$ws2.Range("B80").Value = 74.32655
$ws2.Range("B81").Value = 145.43184
$ws2.Range("B82").Value = 154.2
$ws2.Range("B83").Value = 299.75
$ws2.Range("B86").Value = 12314.43552
$ws2.Range("B87").Value = 10372.09934
$ws2.Range("B88").Value = 10486.69227
$ws2.Range("B89").Value = 1150.9482
$ws2.Range("B90").Value = 1108.46897
$ws2.Range("B91").Value = 284.8668
$ws2.Range("B92").Value = 271.96131
$ws2.Range("B93").Value = 299.98
$ws2.Range("B94").Value = 157.15365
$ws2.Range("B95").Value = 138.42
$ws2.Range("B96").Value = 108.89
$ws2.Range("B97").Value = 105.79
